# Adding "Not applicable" (-1) row to the dbo_sampmet lookup table and
# shifting all existing rows down by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at row 2 (pushes the former row 2..26 down to 3..27)
$ws.Range("A2").EntireRow.Insert()

# Populate the new row 2 with the "Not applicable" lookup entry
$ws.Range("A2").Value = -1
$ws.Range("B2").Value = "Not applicable"

# Keep the named range / used-range in sync with the new last row (27)
$wb.Names.Item("dbo_sampmet").RefersTo = "=dbo_sampmet!`$A`$1:`$G`$27"

# Match the author's final selection on the sheet
[void]$ws.Range("B6").Select()
